$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "index" column and the obsolete leading numeric column (both
# were left over from a stale export); this shifts Recall_data / Recall_jats /
# Recall_multi_lang (and their values) left into columns A:C.
$ws.Range("A:B").Delete()

# Header row: swap the thin box border for a single medium rule under the text.
$hdr = $ws.Range("A1:C1")
$hdr.Borders.LineStyle = 0
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Weight = -4138

# Data row: display the recall values as fixed 4-decimal numbers.
$ws.Range("A2:C2").NumberFormat = "0.0000"

$ws.Range("A1:C2").Select() | Out-Null
